$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7-9 (shift up). This slides the old row10 (which carries the
# unusual "style 3" formatting on column B) up to row 7, old row11 -> row8,
# old row12 -> row9, old row13 -> row10, preserving their existing column-B
# formatting (style 2) exactly as needed by the target layout.
$ws.Range("A7:E9").Delete(-4162)

# Row 2
$ws.Range("A2").Value = "PMS_TC_01"
$ws.Range("B2").Value = "1.5.2"
$ws.Range("C2").Value = "employee_login()"
$ws.Range("D2").Value = "Test_case 1"

# Row 3
$ws.Range("A3").Value = "PMS_TC_02"
$ws.Range("B3").Value = "1.5.2"
$ws.Range("C3").Value = "manager_login()"
$ws.Range("D3").Value = "Test_case 2"

# Row 4
$ws.Range("A4").Value = "PMS_TC_03"
$ws.Range("B4").Value = "             1.5.3.1"
$ws.Range("C4").Value = "employee_registration()"
$ws.Range("D4").Value = "Test_case 1"

# Row 5
$ws.Range("A5").Value = "PMS_TC_04"
$ws.Range("B5").Value = "1.5.4.1"
$ws.Range("C5").Value = "manager_registration()"
$ws.Range("D5").Value = "Test_case 2"

# Row 6 (used to have its 5th column filled in E; now it is D instead)
$ws.Range("A6").Value = "PMS_TC_05"
$ws.Range("B6").Value = "1.5.3.3"
$ws.Range("C6").Value = "Verify the employees ts added to the project"
$ws.Range("E6").Value = ""
$ws.Range("D6").Value = "Test_case 3"

# Row 7 (landed here from the deleted old row 10 - keep its style 3 intact,
# only change the cell values)
$ws.Range("A7").Value = "PMS_TC_06"
$ws.Range("B7").Value = "             1.5.4.4"
$ws.Range("C7").Value = "maintain_project_status()"
$ws.Range("D7").Value = "Test_case 3"

# Row 8
$ws.Range("A8").Value = "PMS_TC_07"
$ws.Range("B8").Value = "1.5.4.5"
$ws.Range("C8").Value = "adding_report()"
$ws.Range("D8").Value = "Test_case 3"

# Row 9 (used to have column D filled; now uses E instead)
$ws.Range("A9").Value = "PMS_TC_08"
$ws.Range("B9").Value = "1.5.2"
$ws.Range("C9").Value = "employee_menu()"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "IT_case1"

# Row 10 (used to have column D filled; now uses E instead)
$ws.Range("A10").Value = "PMS_TC_09"
$ws.Range("B10").Value = "1.5.2"
$ws.Range("C10").Value = "manager_menu()"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = "IT_case 2"

# Rows 11-13 become blank apart from a centred, empty cell in column B
$ws.Range("A11").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B11").HorizontalAlignment = -4108
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""

$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""

$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""

$ws.Range("A11").Select()
